$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 2512.4482
$ws_ALC.Range("J17").Value = 2512.4482
$ws_ALC.Range("L17").Value = 7537.344599999999
$ws_ALC.Range("N17").Value = -7873.344599999999

$ws_ALC.Range("H92").Value = 618.1667
$ws_ALC.Range("I92").Value = 521.9
$ws_ALC.Range("K92").Value = 521.9
$ws_ALC.Range("M92").Value = 726.1

$ws_ALC.Range("H112").Value = 1897.2667
$ws_ALC.Range("J112").Value = 2105.3635
$ws_ALC.Range("L112").Value = 6316.0905
$ws_ALC.Range("N112").Value = -8532.0905

$ws_ALC.Range("H114").Value = 41475.332
$ws_ALC.Range("J114").Value = 41475.332
$ws_ALC.Range("L114").Value = 41475.332
$ws_ALC.Range("N114").Value = -50153.332

$ws_ALC.Range("H117").Value = 48727.8
$ws_ALC.Range("J117").Value = 48727.8
$ws_ALC.Range("L117").Value = 48727.8
$ws_ALC.Range("N117").Value = -57905.8

$ws_ALC.Range("H129").Value = 2171
$ws_ALC.Range("J129").Value = 1734.75
$ws_ALC.Range("L129").Value = 5204.25
$ws_ALC.Range("N129").Value = -15204.25

$ws_ALC.Range("H138").Value = 2247.0588
$ws_ALC.Range("I138").Value = 1510.7142
$ws_ALC.Range("J138").Value = 3436.5386
$ws_ALC.Range("K138").Value = 4532.142599999999
$ws_ALC.Range("L138").Value = 10309.6158
$ws_ALC.Range("M138").Value = 607.8574000000008
$ws_ALC.Range("N138").Value = -20589.6158

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H45").Value = 1715.9642
$ws_ARM.Range("I45").Value = 1536.3158
$ws_ARM.Range("J45").Value = 2095.2222
$ws_ARM.Range("K45").Value = 1536.3158
$ws_ARM.Range("L45").Value = 2095.2222
$ws_ARM.Range("M45").Value = -1159.3158
$ws_ARM.Range("N45").Value = -2849.2222

$ws_ARM.Range("H56").Value = 12000
$ws_ARM.Range("J56").Value = 12000
$ws_ARM.Range("L56").Value = 12000
$ws_ARM.Range("N56").Value = -13484

$ws_ARM.Range("H101").Value = 48054.668
$ws_ARM.Range("J101").Value = 48054.668
$ws_ARM.Range("L101").Value = 48054.668
$ws_ARM.Range("N101").Value = -54544.668

$ws_ARM.Range("H114").Value = 45921.332
$ws_ARM.Range("J114").Value = 45921.332
$ws_ARM.Range("L114").Value = 45921.332
$ws_ARM.Range("N114").Value = -54599.332

$ws_ARM.Range("H121").Value = 34564.832
$ws_ARM.Range("J121").Value = 34564.832
$ws_ARM.Range("L121").Value = 34564.832
$ws_ARM.Range("N121").Value = -38058.832

$ws_ARM.Range("H132").Value = 10418317
$ws_ARM.Range("I132").Value = 19231948
$ws_ARM.Range("K132").Value = 57695844
$ws_ARM.Range("M132").Value = -57693314

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 626.12
$ws_BSM.Range("I94").Value = 582.5238000000001
$ws_BSM.Range("J94").Value = 855
$ws_BSM.Range("K94").Value = 582.5238000000001
$ws_BSM.Range("L94").Value = 855
$ws_BSM.Range("M94").Value = -131.5238000000001
$ws_BSM.Range("N94").Value = -1757

$ws_BSM.Range("H116").Value = 45689.332
$ws_BSM.Range("J116").Value = 45689.332
$ws_BSM.Range("L116").Value = 45689.332
$ws_BSM.Range("N116").Value = -54867.332

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 1747.74
$ws_CRP.Range("I31").Value = 837.2461499999999
$ws_CRP.Range("J31").Value = 3438.6572
$ws_CRP.Range("K31").Value = 837.2461499999999
$ws_CRP.Range("L31").Value = 3438.6572
$ws_CRP.Range("M31").Value = -542.2461499999999
$ws_CRP.Range("N31").Value = -4028.6572

$ws_CRP.Range("H34").Value = 1747.74
$ws_CRP.Range("I34").Value = 837.2461499999999
$ws_CRP.Range("J34").Value = 3438.6572
$ws_CRP.Range("K34").Value = 837.2461499999999
$ws_CRP.Range("L34").Value = 3438.6572
$ws_CRP.Range("M34").Value = -635.2461499999999
$ws_CRP.Range("N34").Value = -3842.6572

$ws_CRP.Range("H80").Value = 33174.6
$ws_CRP.Range("J80").Value = 33174.6
$ws_CRP.Range("L80").Value = 33174.6
$ws_CRP.Range("N80").Value = -35420.6

$ws_CRP.Range("H83").Value = 33174.6
$ws_CRP.Range("J83").Value = 33174.6
$ws_CRP.Range("L83").Value = 99523.79999999999
$ws_CRP.Range("N83").Value = -110755.8

$ws_CRP.Range("H110").Value = 45003.5
$ws_CRP.Range("J110").Value = 45003.5
$ws_CRP.Range("L110").Value = 45003.5
$ws_CRP.Range("N110").Value = -53183.5

$ws_CRP.Range("H122").Value = 189319.72
$ws_CRP.Range("I122").Value = 220754
$ws_CRP.Range("J122").Value = 714
$ws_CRP.Range("K122").Value = 662262
$ws_CRP.Range("L122").Value = 2142
$ws_CRP.Range("M122").Value = -659812
$ws_CRP.Range("N122").Value = -7042

$ws_CRP.Range("H132").Value = 53996.965
$ws_CRP.Range("I132").Value = 1741.9546
$ws_CRP.Range("K132").Value = 5225.8638
$ws_CRP.Range("M132").Value = -2695.8638

$ws_CRP.Range("H137").Value = 68898.57000000001
$ws_CRP.Range("J137").Value = 68898.57000000001
$ws_CRP.Range("L137").Value = 68898.57000000001
$ws_CRP.Range("N137").Value = -79098.57000000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H75").Value = 3420.6667
$ws_CUL.Range("I75").Value = 300
$ws_CUL.Range("J75").Value = 3643.5715
$ws_CUL.Range("K75").Value = 900
$ws_CUL.Range("L75").Value = 10930.7145
$ws_CUL.Range("M75").Value = 98
$ws_CUL.Range("N75").Value = -12926.7145

$ws_CUL.Range("H78").Value = 3420.6667
$ws_CUL.Range("I78").Value = 300
$ws_CUL.Range("J78").Value = 3643.5715
$ws_CUL.Range("K78").Value = 2700
$ws_CUL.Range("L78").Value = 32792.1435
$ws_CUL.Range("M78").Value = 2292
$ws_CUL.Range("N78").Value = -42776.1435

$ws_CUL.Range("H113").Value = 3584.9119
$ws_CUL.Range("I113").Value = 5379.143
$ws_CUL.Range("J113").Value = 686.53845
$ws_CUL.Range("K113").Value = 16137.429
$ws_CUL.Range("L113").Value = 2059.61535
$ws_CUL.Range("M113").Value = -13967.429
$ws_CUL.Range("N113").Value = -6399.61535

$ws_CUL.Range("H131").Value = 2999.8909
$ws_CUL.Range("I131").Value = 33733
$ws_CUL.Range("J131").Value = 1226.8269
$ws_CUL.Range("K131").Value = 101199
$ws_CUL.Range("L131").Value = 3680.4807
$ws_CUL.Range("M131").Value = -96159
$ws_CUL.Range("N131").Value = -13760.4807

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H26").Value = 19907.5
$ws_GSM.Range("J26").Value = 19907.5
$ws_GSM.Range("L26").Value = 19907.5
$ws_GSM.Range("N26").Value = -20467.5

$ws_GSM.Range("H50").Value = 19907.5
$ws_GSM.Range("J50").Value = 19907.5
$ws_GSM.Range("L50").Value = 19907.5
$ws_GSM.Range("N50").Value = -20903.5

$ws_GSM.Range("H113").Value = 986.6667
$ws_GSM.Range("I113").Value = 973.5
$ws_GSM.Range("J113").Value = 1013
$ws_GSM.Range("K113").Value = 973.5
$ws_GSM.Range("L113").Value = 1013
$ws_GSM.Range("M113").Value = 1196.5
$ws_GSM.Range("N113").Value = -5353

$ws_GSM.Range("H114").Value = 37454.832
$ws_GSM.Range("J114").Value = 37454.832
$ws_GSM.Range("L114").Value = 37454.832
$ws_GSM.Range("N114").Value = -46132.832

$ws_GSM.Range("H133").Value = 49000
$ws_GSM.Range("J133").Value = 49000
$ws_GSM.Range("L133").Value = 49000
$ws_GSM.Range("N133").Value = -59120

$ws_GSM.Range("H139").Value = 34031
$ws_GSM.Range("J139").Value = 34031
$ws_GSM.Range("L139").Value = 34031
$ws_GSM.Range("N139").Value = -44311

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H93").Value = 1389.9
$ws_LTW.Range("I93").Value = 1066.75
$ws_LTW.Range("J93").Value = 1507.409
$ws_LTW.Range("K93").Value = 1066.75
$ws_LTW.Range("L93").Value = 1507.409
$ws_LTW.Range("M93").Value = 181.25
$ws_LTW.Range("N93").Value = -4003.409

$ws_LTW.Range("H100").Value = 1515.0834
$ws_LTW.Range("I100").Value = 1515.0834
$ws_LTW.Range("J100").Value = 0
$ws_LTW.Range("K100").Value = 1515.0834
$ws_LTW.Range("L100").Value = 0
$ws_LTW.Range("M100").Value = -974.0834
$ws_LTW.Range("N100").ClearContents()

$ws_LTW.Range("H109").Value = 35273
$ws_LTW.Range("J109").Value = 35273
$ws_LTW.Range("L109").Value = 35273
$ws_LTW.Range("N109").Value = -38047

$ws_LTW.Range("H119").Value = 46728
$ws_LTW.Range("J119").Value = 46728
$ws_LTW.Range("L119").Value = 46728
$ws_LTW.Range("N119").Value = -56404

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H86").Value = 30250
$ws_WVR.Range("J86").Value = 30250
$ws_WVR.Range("L86").Value = 30250
$ws_WVR.Range("N86").Value = -32496

$ws_WVR.Range("H89").Value = 30250
$ws_WVR.Range("J89").Value = 30250
$ws_WVR.Range("L89").Value = 151250
$ws_WVR.Range("N89").Value = -162482

$ws_WVR.Range("H117").Value = 39958.4
$ws_WVR.Range("J117").Value = 39958.4
$ws_WVR.Range("L117").Value = 39958.4
$ws_WVR.Range("N117").Value = -49136.4

$ws_WVR.Range("H138").Value = 42473.75
$ws_WVR.Range("J138").Value = 42473.75
$ws_WVR.Range("L138").Value = 42473.75
$ws_WVR.Range("N138").Value = -52753.75

$ws_WVR.Range("H139").Value = 57650
$ws_WVR.Range("J139").Value = 57650
$ws_WVR.Range("L139").Value = 57650
$ws_WVR.Range("N139").Value = -67930
